$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1079
$ws.Range("I38").Value = 255.71428
$ws.Range("K38").Value = 767.14284
$ws.Range("M38").Value = -395.14284
$ws.Range("H62").Value = 3657.5
$ws.Range("I62").Value = 3293.3333
$ws.Range("K62").Value = 3293.3333
$ws.Range("M62").Value = -2669.3333
$ws.Range("H65").Value = 3657.5
$ws.Range("I65").Value = 3293.3333
$ws.Range("K65").Value = 16466.6665
$ws.Range("M65").Value = -13346.6665
$ws.Range("H95").Value = 35828.125
$ws.Range("J95").Value = 35828.125
$ws.Range("L95").Value = 35828.125
$ws.Range("N95").Value = -41320.125
$ws.Range("H98").Value = 1559.7307
$ws.Range("I98").Value = 1254.9048
$ws.Range("J98").Value = 2840
$ws.Range("K98").Value = 1254.9048
$ws.Range("L98").Value = 2840
$ws.Range("M98").Value = 243.0952
$ws.Range("N98").Value = -5836
$ws.Range("H106").Value = 536.875
$ws.Range("I106").Value = 327.85715
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 327.85715
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = 303.14285
$ws.Range("N106").Value = -3262
$ws.Range("H122").Value = 1559.7307
$ws.Range("I122").Value = 1254.9048
$ws.Range("J122").Value = 2840
$ws.Range("K122").Value = 3764.7144
$ws.Range("L122").Value = 8520
$ws.Range("M122").Value = -1314.7144
$ws.Range("N122").Value = -13420
$ws.Range("H129").Value = 3290576.8
$ws.Range("I129").Value = 35715332
$ws.Range("J129").Value = 1109.0435
$ws.Range("K129").Value = 107145996
$ws.Range("L129").Value = 3327.1305
$ws.Range("M129").Value = -107140996
$ws.Range("N129").Value = -13327.1305
$ws.Range("H132").Value = 5131180.5
$ws.Range("I132").Value = 6063344.5
$ws.Range("K132").Value = 18190033.5
$ws.Range("M132").Value = -18187503.5
$ws.Range("H137").Value = 7151100
$ws.Range("I137").Value = 10008950
$ws.Range("J137").Value = 6474
$ws.Range("K137").Value = 30026850
$ws.Range("L137").Value = 19422
$ws.Range("M137").Value = -30024300
$ws.Range("N137").Value = -24522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7527.506
$ws.Range("I32").Value = 5116.273
$ws.Range("K32").Value = 5116.273
$ws.Range("M32").Value = -4829.273
$ws.Range("H45").Value = 1453.3948
$ws.Range("I45").Value = 1104.1111
$ws.Range("J45").Value = 2310.7273
$ws.Range("K45").Value = 1104.1111
$ws.Range("L45").Value = 2310.7273
$ws.Range("M45").Value = -727.1111000000001
$ws.Range("N45").Value = -3064.7273
$ws.Range("H61").Value = 2999.111
$ws.Range("I61").Value = 1517
$ws.Range("J61").Value = 5963.3335
$ws.Range("K61").Value = 1517
$ws.Range("L61").Value = 5963.3335
$ws.Range("M61").Value = -1305
$ws.Range("N61").Value = -6387.3335
$ws.Range("H63").Value = 2911.647
$ws.Range("I63").Value = 2083.1667
$ws.Range("K63").Value = 2083.1667
$ws.Range("M63").Value = -1397.1667
$ws.Range("H66").Value = 2911.647
$ws.Range("I66").Value = 2083.1667
$ws.Range("K66").Value = 10415.8335
$ws.Range("M66").Value = -6983.833500000001
$ws.Range("H97").Value = 551.48
$ws.Range("I97").Value = 440.4
$ws.Range("K97").Value = 440.4
$ws.Range("M97").Value = 55.60000000000002
$ws.Range("H136").Value = 2999.111
$ws.Range("I136").Value = 1517
$ws.Range("J136").Value = 5963.3335
$ws.Range("K136").Value = 4551
$ws.Range("L136").Value = 17890.0005
$ws.Range("M136").Value = -2001
$ws.Range("N136").Value = -22990.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 458.56522
$ws.Range("I94").Value = 492.9375
$ws.Range("J94").Value = 380
$ws.Range("K94").Value = 492.9375
$ws.Range("L94").Value = 380
$ws.Range("M94").Value = -41.9375
$ws.Range("N94").Value = -1282
$ws.Range("H105").Value = 1883.0769
$ws.Range("I105").Value = 1698.1818
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 1698.1818
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = 48.81819999999993
$ws.Range("N105").Value = -6394

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2329177.5
$ws.Range("I31").Value = 4001764.8
$ws.Range("K31").Value = 4001764.8
$ws.Range("M31").Value = -4001469.8
$ws.Range("H34").Value = 2329177.5
$ws.Range("I34").Value = 4001764.8
$ws.Range("K34").Value = 4001764.8
$ws.Range("M34").Value = -4001562.8
$ws.Range("H58").Value = 12198728
$ws.Range("I58").Value = 2281.08
$ws.Range("K58").Value = 2281.08
$ws.Range("M58").Value = -2078.08
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -9594
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -8596
$ws.Range("N91").ClearContents()
$ws.Range("H94").Value = 1454.069
$ws.Range("I94").Value = 1178
$ws.Range("J94").Value = 1526.0869
$ws.Range("K94").Value = 1178
$ws.Range("L94").Value = 1526.0869
$ws.Range("M94").Value = -727
$ws.Range("N94").Value = -2428.0869
$ws.Range("H122").Value = 1371.0975
$ws.Range("I122").Value = 1118.6364
$ws.Range("K122").Value = 3355.9092
$ws.Range("M122").Value = -905.9092000000001
$ws.Range("H132").Value = 2059.2974
$ws.Range("I132").Value = 1853.1818
$ws.Range("J132").Value = 2361.6
$ws.Range("K132").Value = 5559.5454
$ws.Range("L132").Value = 7084.799999999999
$ws.Range("M132").Value = -3029.5454
$ws.Range("N132").Value = -12144.8
$ws.Range("H134").Value = 2503.0454
$ws.Range("I134").Value = 1386.2941
$ws.Range("K134").Value = 4158.8823
$ws.Range("M134").Value = -1623.8823
$ws.Range("H136").Value = 12198728
$ws.Range("I136").Value = 2281.08
$ws.Range("K136").Value = 6843.24
$ws.Range("M136").Value = -4293.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 214.47368
$ws.Range("I38").Value = 84
$ws.Range("J38").Value = 261.07144
$ws.Range("K38").Value = 252
$ws.Range("L38").Value = 783.21432
$ws.Range("M38").Value = 95
$ws.Range("N38").Value = -1477.21432
$ws.Range("H68").Value = 6800.4
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6800.4
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 20401.2
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -22023.2
$ws.Range("H71").Value = 6800.4
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6800.4
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 61203.6
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -69315.60000000001
$ws.Range("H122").Value = 1537.5385
$ws.Range("I122").Value = 930
$ws.Range("J122").Value = 1719.8
$ws.Range("K122").Value = 8370
$ws.Range("L122").Value = 15478.2
$ws.Range("M122").Value = -5920
$ws.Range("N122").Value = -20378.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 32462.545
$ws.Range("I102").Value = 1748.9584
$ws.Range("J102").Value = 114365.445
$ws.Range("K102").Value = 1748.9584
$ws.Range("L102").Value = 114365.445
$ws.Range("M102").Value = -126.9584
$ws.Range("N102").Value = -117609.445
$ws.Range("H107").Value = 1324.5834
$ws.Range("I107").Value = 261.875
$ws.Range("J107").Value = 3450
$ws.Range("K107").Value = 261.875
$ws.Range("L107").Value = 3450
$ws.Range("M107").Value = 1658.125
$ws.Range("N107").Value = -7290
$ws.Range("H113").Value = 1789.2084
$ws.Range("I113").Value = 1067.2858
$ws.Range("J113").Value = 2799.9
$ws.Range("K113").Value = 1067.2858
$ws.Range("L113").Value = 2799.9
$ws.Range("M113").Value = 1102.7142
$ws.Range("N113").Value = -7139.9
$ws.Range("H122").Value = 3949.7932
$ws.Range("I122").Value = 3024
$ws.Range("K122").Value = 9072
$ws.Range("M122").Value = -6622
$ws.Range("H126").Value = 2886.15
$ws.Range("I126").Value = 1203
$ws.Range("J126").Value = 4008.25
$ws.Range("K126").Value = 3609
$ws.Range("L126").Value = 12024.75
$ws.Range("M126").Value = -1139
$ws.Range("N126").Value = -16964.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1892.579
$ws.Range("I7").Value = 1161.6666
$ws.Range("J7").Value = 2550.4
$ws.Range("K7").Value = 1161.6666
$ws.Range("L7").Value = 2550.4
$ws.Range("M7").Value = -1049.6666
$ws.Range("N7").Value = -2774.4
$ws.Range("H40").Value = 2940
$ws.Range("I40").Value = 2233.3333
$ws.Range("K40").Value = 2233.3333
$ws.Range("M40").Value = -2097.3333
$ws.Range("H41").Value = 250000
$ws.Range("I41").Value = 250000
$ws.Range("K41").Value = 250000
$ws.Range("M41").Value = -249562
$ws.Range("H122").Value = 2828.524
$ws.Range("I122").Value = 2371
$ws.Range("K122").Value = 7113
$ws.Range("M122").Value = -4663
$ws.Range("H126").Value = 1892.579
$ws.Range("I126").Value = 1161.6666
$ws.Range("J126").Value = 2550.4
$ws.Range("K126").Value = 3484.9998
$ws.Range("L126").Value = 7651.200000000001
$ws.Range("M126").Value = -1014.9998
$ws.Range("N126").Value = -12591.2
$ws.Range("H132").Value = 2391.879
$ws.Range("I132").Value = 1447.5
$ws.Range("J132").Value = 3525.1333
$ws.Range("K132").Value = 4342.5
$ws.Range("L132").Value = 10575.3999
$ws.Range("M132").Value = -1812.5
$ws.Range("N132").Value = -15635.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 590571.3
$ws.Range("I122").Value = 835202.2
$ws.Range("K122").Value = 2505606.6
$ws.Range("M122").Value = -2503156.6
$ws.Range("H126").Value = 3848186.8
$ws.Range("I126").Value = 1225.9412
$ws.Range("K126").Value = 3677.8236
$ws.Range("M126").Value = -1207.8236
$ws.Range("H132").Value = 236903.98
$ws.Range("I132").Value = 358298.9
$ws.Range("J132").Value = 10300.134
$ws.Range("K132").Value = 1074896.7
$ws.Range("L132").Value = 30900.402
$ws.Range("M132").Value = -1072366.7
$ws.Range("N132").Value = -35960.402
